$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "1.001", "313.64").
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.369.96'
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").Value = '1.824.80'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '313.64'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.4470'
$ws.Range("E7").Value = '  +3.19%  '
$ws.Range("D8").Value = '0.3753'
$ws.Range("E8").Value = '  +2.00%  '
$ws.Range("D9").Value = '0.07507'
$ws.Range("E9").Value = '  +2.92%  '
$ws.Range("D10").Value = '0.8891'
$ws.Range("E10").Value = '  +4.99%  '
$ws.Range("D11").Value = '21.04'
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").Value = '1.824.41'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '6.761'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").Value = '94.01'
$ws.Range("E14").Value = '  +5.08%  '
$ws.Range("D15").Value = '5.413'
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").Value = '0.07102'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '0.000008815'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").Value = '27.367.84'
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").Value = '5.264'
$ws.Range("E22").Value = '  +2.29%  '
$ws.Range("D23").Value = '10.92'
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '2.057.02'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = '1.969'
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("D26").Value = '2.379'
$ws.Range("E26").Value = '  +7.28%  '
$ws.Range("D27").Value = '151.42'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("D29").Value = '5.355'
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("D30").Value = '118.10'
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").Value = '0.08833'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = '0.7854'
$ws.Range("E32").Value = '  +5.91%  '
$ws.Range("D33").Value = '1.199'
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("D34").Value = '4.517'
$ws.Range("E34").Value = '  +1.59%  '
$ws.Range("D35").Value = '2.940'
$ws.Range("E35").Value = '  +1.02%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = '1.109'
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("D38").Value = '0.01995'
$ws.Range("E38").Value = '  +2.29%  '
$ws.Range("D39").Value = '0.05331'
$ws.Range("E39").Value = '  +1.64%  '
$ws.Range("D40").Value = '7.382'
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").Value = '0.5314'
$ws.Range("E41").Value = '  +3.60%  '
$ws.Range("D42").Value = '0.1728'
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("D43").Value = '2.858'
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").Value = '2.298'
$ws.Range("E44").Value = '  +18.39%  '
$ws.Range("D45").Value = '8.762'
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").Value = '0.5112'
$ws.Range("E46").Value = '  +7.27%  '
$ws.Range("D47").Value = '10.68'
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '105.84'
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.701'
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").Value = '0.06378'
$ws.Range("E51").Value = '  +0.62%  '
